$p = $ppt.ActivePresentation
$p.Slides.Add(3, 1)
